# ReportDefs_vervestacks.xlsx - scenario set update
# Switches the ScenMap sheet's scenario-group from the "ngfs" 5-scenario set
# (Postponed Transition / Target Net Zero 2050 / Declared NDCs / Limited to
# 2 deg / Current Policies) to the new "ar6_r10" 5-scenario set
# (e 1.5 deg no OS / d 1.5 deg OS / c 2 deg (67%) / b 2 deg (50%) / a 3 deg),
# and records the long-form IPCC AR6 category description for each scenario
# in column L.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenMap")

# Scenario-group name (drives the "sg_" prefixed labels in I5/J5)
$ws.Range("I2").Value = "ar6_r10"

# Short scenario codes used throughout the sheet (I6:I10 feed I11:I55 via
# formulas, and those in turn feed B/C/H columns, so only these five cells
# need to be edited directly).
$ws.Range("I6").Value = "e 1.5 deg no OS"
$ws.Range("I7").Value = "d 1.5 deg OS"
$ws.Range("I8").Value = "c 2 deg (67%)"
$ws.Range("I9").Value = "b 2 deg (50%)"
$ws.Range("I10").Value = "a 3 deg"

# Long-form AR6 category descriptions (new column L)
$ws.Range("L6").Value = "Limit warming to 1.5°C (>50%) with no or limited overshoot"
$ws.Range("L7").Value = "Limit warming to 1.5°C (>67%) with high overshoot"
$ws.Range("L8").Value = "Limit warming to 2°C (>67%) with higher action post-2030"
$ws.Range("L9").Value = "Limit warming to 2°C (>50%) with immediate action"
$ws.Range("L10").Value = "Likely above 3°C warming with limited mitigation"

# Update the saved cursor/selection position
$ws.Range("A3").Select()
